$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 5623.7055
    3  = 5458.229
    4  = 5507.9745
    5  = 5475.498
    6  = 5469.52
    7  = 5575.703
    8  = 6141.8735
    9  = 7163.6355
    10 = 8647.460500000001
    11 = 12435.7905
    12 = 15471.022
    13 = 16384.8475
    14 = 16062.228
    15 = 16185.204
    16 = 16525.544
    17 = 16462.6175
    18 = 16761.7275
    19 = 16532.1625
    20 = 15890.3885
    21 = 13683.775
    22 = 11966.283
    23 = 9398.494000000001
    24 = 6574.9635
    25 = 5831.934499999999
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
